$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1 - same style as the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats (copy formatting only, not value)
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 13:41:49.749943",
    "2021-10-05 13:41:49.749955",
    "2021-10-05 13:41:49.749959",
    "2021-10-05 13:41:49.749963",
    "2021-10-05 13:41:49.749966",
    "2021-10-05 13:41:49.749969",
    "2021-10-05 13:41:49.749972",
    "2021-10-05 13:41:49.749975",
    "2021-10-05 13:41:49.749978",
    "2021-10-05 13:41:49.749981",
    "2021-10-05 13:41:49.749984",
    "2021-10-05 13:41:49.749987",
    "2021-10-05 13:41:49.749990",
    "2021-10-05 13:41:49.749993",
    "2021-10-05 13:41:49.749996",
    "2021-10-05 13:41:49.749998",
    "2021-10-05 13:41:49.750001",
    "2021-10-05 13:41:49.750005",
    "2021-10-05 13:41:49.750007",
    "2021-10-05 13:41:49.750011",
    "2021-10-05 13:41:49.750013"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
